# Update the "Training Dashboard" sheet with the new progress date (04-Nov-2025).
# For each data row (3 through 25):
#   - column H ("PERIOD TO EXPIRE") is decremented by 1
#   - column I ("LAST UPDATE") is set to "04-Nov-2025"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Training Dashboard")

for ($row = 3; $row -le 25; $row++) {
    $hCell = $ws.Cells.Item($row, 8)   # column H
    $hCell.Value2 = $hCell.Value2 - 1

    $iCell = $ws.Cells.Item($row, 9)   # column I
    $iCell.NumberFormat = "@"
    $iCell.Value2 = "04-Nov-2025"
}
